# Slide 2 ("For deploying a web app on Windows VM (Azure)"), the
# "Content Placeholder 2" shape's first bullet.
#   Before: "Create a " + "Windows based VM"   (2 runs)
#   After:  "Create "   + "a Windows-based VM" (2 runs, same rPr each)
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$para1 = $tr.Paragraphs(1, 1)

# Edit exactly on the existing run boundaries so each run keeps its own
# rPr (and PowerPoint doesn't stamp a fresh dirty="0" on the 2nd run).
$run1 = $para1.Characters(1, 9)
$run1.Text = "Create "

$run2 = $para1.Characters(8, 16)
$run2.Text = "a Windows-based VM"
